$wb = $excel.ActiveWorkbook

# sheet1 - 展览 (Exhibition)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 277
$ws.Range("F4").Value = 627
$ws.Range("F5").Value = 2772
$ws.Range("F7").Value = 216
$ws.Range("F8").Value = 269
$ws.Range("F9").Value = 6350
$ws.Range("F11").Value = 77
$ws.Range("F13").Value = 5019
$ws.Range("F14").Value = 103
$ws.Range("F15").Value = 548
$ws.Range("F16").Value = 2624
$ws.Range("F18").Value = 1516
$ws.Range("F19").Value = 1220
$ws.Range("F20").Value = 316
$ws.Range("F23").Value = 1067
$ws.Range("F24").Value = 240
$ws.Range("F26").Value = 534
$ws.Range("F27").Value = 1377
$ws.Range("F28").Value = 1044
$ws.Range("F29").Value = 2100
$ws.Range("F30").Value = 319
$ws.Range("F31").Value = 578
$ws.Range("F32").Value = 26
$ws.Range("F33").Value = 25
$ws.Range("F34").Value = 89
$ws.Range("F36").Value = 1496
$ws.Range("F42").Value = 296
$ws.Range("F43").Value = 2280
$ws.Range("F44").Value = 2553
$ws.Range("F46").Value = 129
$ws.Range("F47").Value = 273
$ws.Range("F49").Value = 96

# sheet2 - 演出 (Performance)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F7").Value = 405
$ws.Range("F10").Value = 169
$ws.Range("F19").Value = 151
$ws.Range("F20").Value = 38
$ws.Range("F23").Value = 325
$ws.Range("F24").Value = 377
$ws.Range("F26").Value = 9
$ws.Range("F31").Value = 22
$ws.Range("F33").Value = 6

# sheet3 - 本地生活 (Local life)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F6").Value = 1698
$ws.Range("F7").Value = 569
$ws.Range("F8").Value = 1499
$ws.Range("F10").Value = 2520
$ws.Range("F11").Value = 848
$ws.Range("F12").Value = 720
$ws.Range("F13").Value = 24

# sheet4 - 全部类型 (All types)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 277
$ws.Range("F5").Value = 627
$ws.Range("F6").Value = 569
$ws.Range("F7").Value = 2772
$ws.Range("F8").Value = 216
$ws.Range("F9").Value = 1499
$ws.Range("F10").Value = 269
$ws.Range("F11").Value = 2520
$ws.Range("F12").Value = 6350
$ws.Range("F13").Value = 848
$ws.Range("F16").Value = 5019
$ws.Range("F17").Value = 2624
$ws.Range("F19").Value = 1516
$ws.Range("F20").Value = 1220
$ws.Range("F24").Value = 240
$ws.Range("F27").Value = 1377
$ws.Range("F28").Value = 1044
$ws.Range("F29").Value = 2100
$ws.Range("F30").Value = 319
$ws.Range("F31").Value = 578
$ws.Range("F32").Value = 26
$ws.Range("F34").Value = 25
$ws.Range("F41").Value = 325
$ws.Range("F42").Value = 296
$ws.Range("F44").Value = 2280
$ws.Range("F45").Value = 2553
$ws.Range("F46").Value = 129
$ws.Range("F47").Value = 273
$ws.Range("F49").Value = 6
